$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix markdown line breaks in the GCS coverage description strings ---
# Replace literal "<br>" separators with real newlines so the markdown
# renders properly. These strings are duplicated across 12 rows each.

$newLow = "Supports the GCS if coverage is **Low** `nOther members: Global South + EU `n(25-33% of world emissions)"
$newMid = "Supports the GCS if coverage is **Mid** `nGlobal South + China `n(56% of world emissions)"
$newHigh = "Supports the GCS if coverage is **High** `nGlobal South + China + EU + various HICs `n(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$newHighColor = "Supports the GCS if coverage is **High**, **color** variant `nGlobal South + China + EU + various HICs `n+ Distributive effects shown using colors on world map"

for ($r = 26; $r -le 37; $r++) { $ws.Range("E$r").Value = $newLow }
for ($r = 38; $r -le 49; $r++) { $ws.Range("E$r").Value = $newMid }
for ($r = 50; $r -le 61; $r++) { $ws.Range("E$r").Value = $newHigh }
for ($r = 62; $r -le 73; $r++) { $ws.Range("E$r").Value = $newHighColor }

# --- Update re-computed bootstrap summary statistics ---
$ws.Range("B2").Value = 67.7971034621066
$ws.Range("C2").Value = 66.4060461449477
$ws.Range("D2").Value = 69.1881607792655

$ws.Range("B12").Value = 73.8102295238027
$ws.Range("C12").Value = 69.9173212117554
$ws.Range("D12").Value = 77.70313783585

$ws.Range("B14").Value = 55.3296942977469
$ws.Range("C14").Value = 53.8471289764683
$ws.Range("D14").Value = 56.8122596190256

$ws.Range("B24").Value = 48.5761024007792
$ws.Range("C24").Value = 44.2083508526147
$ws.Range("D24").Value = 52.9438539489438

$ws.Range("B38").Value = 67.1068857857389
$ws.Range("C38").Value = 65.398148692157
$ws.Range("D38").Value = 68.8156228793208

$ws.Range("B48").Value = 63.180865855067
$ws.Range("C48").Value = 57.4064206273406
$ws.Range("D48").Value = 68.9553110827935

$ws.Range("B50").Value = 68.4640848889378
$ws.Range("C50").Value = 66.8324549522336
$ws.Range("D50").Value = 70.095714825642

$ws.Range("B60").Value = 59.8978761247204
$ws.Range("C60").Value = 53.8467358677426
$ws.Range("D60").Value = 65.9490163816982

$ws.Range("B62").Value = 61.8966560897
$ws.Range("C62").Value = 60.1646532127507
$ws.Range("D62").Value = 63.6286589666493

$ws.Range("B72").Value = 54.1378328703989
$ws.Range("C72").Value = 47.9978727513762
$ws.Range("D72").Value = 60.2777929894216
